# Applies the cryptos.xlsx data refresh described by the commit:
#   "Updated cryptos list on Fri Sep 29 20:29:50 UTC 2023 with GitHub Actions"
#
# Updates Price (column D) and Volume(1h) (column E) values for the coins
# that changed, including the FraxShare/PaxDollar row swap (rows 40-41).
#
# Cells whose new value looks like a plain number ("215.49", "0.0496", ...)
# are forced to Text format first so Excel keeps storing them as the exact
# original string (matching the source t="inlineStr" cells) instead of
# silently converting them to floating point numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.907.74'
$ws.Range("E2").Value = '  -0.80%  '
$ws.Range("D3").Value = '1.666.19'
$ws.Range("E3").Value = '  +0.48%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.49'
$ws.Range("E5").Value = '  +0.10%  '
$ws.Range("E6").Value = '  +4.90%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("E8").Value = '  +0.93%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.248'
$ws.Range("E9").Value = '  -0.54%  '
$ws.Range("E10").Value = '  +2.64%  '
$ws.Range("E11").Value = '  +3.73%  '
$ws.Range("D12").Value = '1.900.88'
$ws.Range("E12").Value = '  +0.48%  '
$ws.Range("D13").Value = '1.690.16'
$ws.Range("E13").Value = '  +1.91%  '
$ws.Range("E14").Value = '  -0.15%  '
$ws.Range("E15").Value = '  +0.74%  '
$ws.Range("E16").Value = '  +1.37%  '
$ws.Range("D17").Value = '26.907.60'
$ws.Range("E17").Value = '  -0.71%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '234.44'
$ws.Range("E18").Value = '  -1.68%  '
$ws.Range("E19").Value = '  +1.03%  '
$ws.Range("E20").Value = '  +0.37%  '
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.35'
$ws.Range("E22").Value = '  -2.03%  '
$ws.Range("E23").Value = '  -2.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.11'
$ws.Range("E24").Value = '  -1.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.95'
$ws.Range("E25").Value = '  +0.27%  '
$ws.Range("E26").Value = '  -0.40%  '
$ws.Range("E27").Value = '  +1.16%  '
$ws.Range("E28").Value = '  +0.15%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0496'
$ws.Range("E30").Value = '  -0.20%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.17'
$ws.Range("E31").Value = '  +0.11%  '
$ws.Range("D33").Value = '1.455.94'
$ws.Range("E33").Value = '  -4.46%  '
$ws.Range("E34").Value = '  +2.61%  '
$ws.Range("E35").Value = '  +1.52%  '
$ws.Range("E36").Value = '  -0.16%  '
$ws.Range("E37").Value = '  +0.25%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.902'
$ws.Range("E38").Value = '  +1.39%  '
$ws.Range("E39").Value = '  +0.00%  '
$ws.Range("B40").Value = 'PaxDollar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  +0.07%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.71'
$ws.Range("E41").Value = '  -4.07%  '
$ws.Range("E42").Value = '  +0.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.975'
$ws.Range("E43").Value = '  +6.24%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '65.79'
$ws.Range("E44").Value = '  -0.68%  '
$ws.Range("D45").Value = '1.808.63'
$ws.Range("E45").Value = '  +0.55%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.781'
$ws.Range("E46").Value = '  +0.18%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.56'
$ws.Range("E47").Value = '  +0.59%  '
$ws.Range("E48").Value = '  +0.51%  '
$ws.Range("E49").Value = '  -1.19%  '
$ws.Range("E50").Value = '  +4.57%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0505'
$ws.Range("E51").Value = '  +0.01%  '
